{"js": "// Update the worksheet date heading and every arithmetic-problem cell in\n// the single table with the new values from the target revision.\n// The table is read row-by-row, left-to-right (20 rows x 5 columns = 100\n// cells), matching Word.Table#values ordering.\n\nconst NEW_ROWS = [\n  [\"52+13=\", \"30+47=\", \"41+15=\", \"27+35=\", \"68-52=\"],\n  [\"58+1=\", \"9+57=\", \"20-8=\", \"80-37=\", \"19-18=\"],\n  [\"47-3=\", \"9+2=\", \"34+57=\", \"18+35=\", \"18+19=\"],\n  [\"90-51=\", \"76-20=\", \"27+5=\", \"65+21=\", \"74-71=\"],\n  [\"36+29=\", \"15+23=\", \"42-9=\", \"87-10=\", \"75-20=\"],\n  [\"18+36=\", \"91-79=\", \"80-78=\", \"20+15=\", \"68-0=\"],\n  [\"44+7=\", \"3+6=\", \"62+24=\", \"65+29=\", \"4+87=\"],\n  [\"91-24=\", \"95-73=\", \"92-87=\", \"98-46=\", \"42+15=\"],\n  [\"26+6=\", \"4+3=\", \"71-20=\", \"36+28=\", \"92+0=\"],\n  [\"87+2=\", \"26+52=\", \"68-16=\", \"76-25=\", \"74-40=\"],\n  [\"63-37=\", \"77-5=\", \"98-97=\", \"5+55=\", \"33-10=\"],\n  [\"9+88=\", \"54+3=\", \"18+47=\", \"70-2=\", \"49-15=\"],\n  [\"32+58=\", \"75-15=\", \"99-49=\", \"17+32=\", \"0+80=\"],\n  [\"11+39=\", \"79+9=\", \"77-63=\", \"89-73=\", \"91-19=\"],\n  [\"68-63=\", \"54+43=\", \"18+74=\", \"68+30=\", \"72+3=\"],\n  [\"57-4=\", \"13+2=\", \"57+18=\", \"61-51=\", \"7+21=\"],\n  [\"69-27=\", \"15-5=\", \"46+13=\", \"95-42=\", \"22+34=\"],\n  [\"60+26=\", \"75+20=\", \"71-23=\", \"10+7=\", \"70+19=\"],\n  [\"68-26=\", \"78+6=\", \"32+1=\", \"17+58=\", \"75-47=\"],\n  [\"13+30=\", \"86-17=\", \"59-24=\", \"65+3=\", \"55+43=\"],\n];\n\n// 1. Update the date heading paragraph (first paragraph of the body).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.load(\"text\");\nawait context.sync();\n\nif (dateParagraph.text === \"2025-04-03 Thursday\") {\n  dateParagraph.insertText(\"2025-04-04 Friday\", Word.InsertLocation.replace);\n}\n\n// 2. Update every cell of the (only) table with the new arithmetic\n// expressions, preserving all existing formatting (fonts, size, etc.)\n// since only the cell text content is being assigned.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.values = NEW_ROWS;\n\nawait context.sync();\n", "ps1": "# Update the worksheet date heading and every arithmetic-problem cell in\n# the single table with the new values from the target revision.\n# Every old value is unique within the document, so a simple whole-document\n# Find/Replace per pair unambiguously targets the right run while leaving\n# all formatting (fonts, size, paragraph/table structure) untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ old = \"2025-04-03 Thursday\"; new = \"2025-04-04 Friday\" }\n    @{ old = \"7+58=\"; new = \"52+13=\" }\n    @{ old = \"34-30=\"; new = \"30+47=\" }\n    @{ old = \"48+6=\"; new = \"41+15=\" }\n    @{ old = \"5+44=\"; new = \"27+35=\" }\n    @{ old = \"69+3=\"; new = \"68-52=\" }\n    @{ old = \"43+52=\"; new = \"58+1=\" }\n    @{ old = \"5+37=\"; new = \"9+57=\" }\n    @{ old = \"94-73=\"; new = \"20-8=\" }\n    @{ old = \"51-5=\"; new = \"80-37=\" }\n    @{ old = \"8+42=\"; new = \"19-18=\" }\n    @{ old = \"29+24=\"; new = \"47-3=\" }\n    @{ old = \"22+66=\"; new = \"9+2=\" }\n    @{ old = \"83-32=\"; new = \"34+57=\" }\n    @{ old = \"73-0=\"; new = \"18+35=\" }\n    @{ old = \"50-50=\"; new = \"18+19=\" }\n    @{ old = \"98-84=\"; new = \"90-51=\" }\n    @{ old = \"90-8=\"; new = \"76-20=\" }\n    @{ old = \"16-8=\"; new = \"27+5=\" }\n    @{ old = \"24+5=\"; new = \"65+21=\" }\n    @{ old = \"54+18=\"; new = \"74-71=\" }\n    @{ old = \"50-6=\"; new = \"36+29=\" }\n    @{ old = \"64+22=\"; new = \"15+23=\" }\n    @{ old = \"25+60=\"; new = \"42-9=\" }\n    @{ old = \"29+56=\"; new = \"87-10=\" }\n    @{ old = \"64-26=\"; new = \"75-20=\" }\n    @{ old = \"40+18=\"; new = \"18+36=\" }\n    @{ old = \"66-2=\"; new = \"91-79=\" }\n    @{ old = \"90-24=\"; new = \"80-78=\" }\n    @{ old = \"85-19=\"; new = \"20+15=\" }\n    @{ old = \"29+8=\"; new = \"68-0=\" }\n    @{ old = \"19-11=\"; new = \"44+7=\" }\n    @{ old = \"48+0=\"; new = \"3+6=\" }\n    @{ old = \"40+9=\"; new = \"62+24=\" }\n    @{ old = \"98-11=\"; new = \"65+29=\" }\n    @{ old = \"19+36=\"; new = \"4+87=\" }\n    @{ old = \"80-25=\"; new = \"91-24=\" }\n    @{ old = \"61+11=\"; new = \"95-73=\" }\n    @{ old = \"98-81=\"; new = \"92-87=\" }\n    @{ old = \"13+71=\"; new = \"98-46=\" }\n    @{ old = \"77-21=\"; new = \"42+15=\" }\n    @{ old = \"24+33=\"; new = \"26+6=\" }\n    @{ old = \"24-0=\"; new = \"4+3=\" }\n    @{ old = \"11-3=\"; new = \"71-20=\" }\n    @{ old = \"0+52=\"; new = \"36+28=\" }\n    @{ old = \"66-65=\"; new = \"92+0=\" }\n    @{ old = \"94-48=\"; new = \"87+2=\" }\n    @{ old = \"33+52=\"; new = \"26+52=\" }\n    @{ old = \"73+2=\"; new = \"68-16=\" }\n    @{ old = \"95+3=\"; new = \"76-25=\" }\n    @{ old = \"10-7=\"; new = \"74-40=\" }\n    @{ old = \"43+50=\"; new = \"63-37=\" }\n    @{ old = \"59-15=\"; new = \"77-5=\" }\n    @{ old = \"21+29=\"; new = \"98-97=\" }\n    @{ old = \"42+30=\"; new = \"5+55=\" }\n    @{ old = \"21+39=\"; new = \"33-10=\" }\n    @{ old = \"31+23=\"; new = \"9+88=\" }\n    @{ old = \"89-85=\"; new = \"54+3=\" }\n    @{ old = \"44+54=\"; new = \"18+47=\" }\n    @{ old = \"39-2=\"; new = \"70-2=\" }\n    @{ old = \"86-22=\"; new = \"49-15=\" }\n    @{ old = \"37+29=\"; new = \"32+58=\" }\n    @{ old = \"44+46=\"; new = \"75-15=\" }\n    @{ old = \"25-19=\"; new = \"99-49=\" }\n    @{ old = \"82-63=\"; new = \"17+32=\" }\n    @{ old = \"84-74=\"; new = \"0+80=\" }\n    @{ old = \"28+45=\"; new = \"11+39=\" }\n    @{ old = \"6+14=\"; new = \"79+9=\" }\n    @{ old = \"50+21=\"; new = \"77-63=\" }\n    @{ old = \"26+1=\"; new = \"89-73=\" }\n    @{ old = \"24+48=\"; new = \"91-19=\" }\n    @{ old = \"92-88=\"; new = \"68-63=\" }\n    @{ old = \"84-81=\"; new = \"54+43=\" }\n    @{ old = \"56-16=\"; new = \"18+74=\" }\n    @{ old = \"61-40=\"; new = \"68+30=\" }\n    @{ old = \"9+1=\"; new = \"72+3=\" }\n    @{ old = \"12+58=\"; new = \"57-4=\" }\n    @{ old = \"90-46=\"; new = \"13+2=\" }\n    @{ old = \"52-46=\"; new = \"57+18=\" }\n    @{ old = \"67-2=\"; new = \"61-51=\" }\n    @{ old = \"78-64=\"; new = \"7+21=\" }\n    @{ old = \"52+15=\"; new = \"69-27=\" }\n    @{ old = \"40+44=\"; new = \"15-5=\" }\n    @{ old = \"37-22=\"; new = \"46+13=\" }\n    @{ old = \"89-80=\"; new = \"95-42=\" }\n    @{ old = \"7+78=\"; new = \"22+34=\" }\n    @{ old = \"0+96=\"; new = \"60+26=\" }\n    @{ old = \"28+61=\"; new = \"75+20=\" }\n    @{ old = \"26+25=\"; new = \"71-23=\" }\n    @{ old = \"89-33=\"; new = \"10+7=\" }\n    @{ old = \"16+12=\"; new = \"70+19=\" }\n    @{ old = \"77+3=\"; new = \"68-26=\" }\n    @{ old = \"7+11=\"; new = \"78+6=\" }\n    @{ old = \"5+40=\"; new = \"32+1=\" }\n    @{ old = \"16+53=\"; new = \"17+58=\" }\n    @{ old = \"71-69=\"; new = \"75-47=\" }\n    @{ old = \"34+49=\"; new = \"13+30=\" }\n    @{ old = \"19+51=\"; new = \"86-17=\" }\n    @{ old = \"41-15=\"; new = \"59-24=\" }\n    @{ old = \"20+19=\"; new = \"65+3=\" }\n    @{ old = \"21+31=\"; new = \"55+43=\" }\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.old\n    $find.Replacement.Text = $pair.new\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($pair.old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.new, 2) | Out-Null\n}\n\nWrite-Output \"done\"\n"}
